# Generate Report for Handoff
#
# Marks the "Ready for handoff" rows (7, 8, 9, 11, 12, 14) on both the
# zh-cn and de-de sheets as handed off ("ht" priority/handoff-type), and
# refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps on the Overview sheet and on each locale sheet.

$wb = $excel.ActiveWorkbook

$rows = 7,8,9,11,12,14

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for rows that
# were just (re)handed off.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-29 16:23:14"
}

# zh-cn sheet: mark Priority ("ht") and refresh Latest Handoff Datetime.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-29 16:23:06"
}

# de-de sheet: mark Priority ("ht") and refresh Latest Handoff Datetime.
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-29 16:23:14"
}
